$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4165500938170587
$ws.Range("C2").Value = 0.4532998755493539
$ws.Range("D2").Value = 0.8397313451970363
$ws.Range("E2").Value = 0.4349759063693696
$ws.Range("F2").Value = 0.4187560232757783
$ws.Range("G2").Value = 0.5292044909545227
$ws.Range("H2").Value = 0.4360621853365316
$ws.Range("B3").Value = 0.4425011968473371
$ws.Range("C3").Value = 0.4728680731960298
$ws.Range("D3").Value = 0.7676640602148203
$ws.Range("E3").Value = 0.4495895168240502
$ws.Range("F3").Value = 0.4435718024875882
$ws.Range("G3").Value = 0.5348675326294656
$ws.Range("H3").Value = 0.450039492624489
$ws.Range("B4").Value = 0.3717759532087446
$ws.Range("C4").Value = 0.4930149867371041
$ws.Range("D4").Value = 0.4743663168598267
$ws.Range("E4").Value = 0.4126039193304918
$ws.Range("F4").Value = 0.3750959194988736
$ws.Range("G4").Value = 0.4753858453573038
$ws.Range("H4").Value = 0.4139473131987446
$ws.Range("B5").Value = 0.4802277040075707
$ws.Range("C5").Value = 0.5775338682106308
$ws.Range("D5").Value = 0.364676707402917
$ws.Range("E5").Value = 0.4740925553094573
$ws.Range("F5").Value = 0.479423323321474
$ws.Range("G5").Value = 0.5339142996764678
$ws.Range("H5").Value = 0.473703034265991
$ws.Range("B6").Value = 0.4707275492917793
$ws.Range("C6").Value = 0.5646114238160951
$ws.Range("D6").Value = 0.2923816395688805
$ws.Range("E6").Value = 0.4692799780971428
$ws.Range("F6").Value = 0.4703025743668272
$ws.Range("G6").Value = 0.5062968331428027
$ws.Range("H6").Value = 0.4689871014179974
$ws.Range("B7").Value = 0.4945077918741457
$ws.Range("C7").Value = 0.6222607695612322
$ws.Range("D7").Value = 0.328609636718072
$ws.Range("E7").Value = 0.481179764654688
$ws.Range("F7").Value = 0.4933860710767114
$ws.Range("G7").Value = 0.5625672726102521
$ws.Range("H7").Value = 0.4807044994505518
$ws.Range("B8").Value = 0.4303344470155975
$ws.Range("C8").Value = 0.2492231869568524
$ws.Range("D8").Value = 0.1935278923803382
$ws.Range("E8").Value = 0.4171685643628851
$ws.Range("F8").Value = 0.4291178983441817
$ws.Range("G8").Value = 0.2571093819663559
$ws.Range("H8").Value = 0.4166872841406259
$ws.Range("B9").Value = 0.4982263498511468
$ws.Range("C9").Value = 0.5286358997867628
$ws.Range("D9").Value = 0.5673291140965486
$ws.Range("E9").Value = 0.4756995253731762
$ws.Range("F9").Value = 0.4972499334998906
$ws.Range("G9").Value = 0.4971649692636496
$ws.Range("H9").Value = 0.4753935026063243
$ws.Range("B10").Value = 0.5123575806189912
$ws.Range("C10").Value = 0.5960745203152207
$ws.Range("D10").Value = 0.4268077567395996
$ws.Range("E10").Value = 0.4850680614742685
$ws.Range("F10").Value = 0.5111134471564462
$ws.Range("G10").Value = 0.5948536767995212
$ws.Range("H10").Value = 0.4848457968482907
$ws.Range("B11").Value = 0.5075664387102803
$ws.Range("C11").Value = 0.5990363700050438
$ws.Range("D11").Value = 0.3745725955835346
$ws.Range("E11").Value = 0.4846872260551944
$ws.Range("F11").Value = 0.5063591113775412
$ws.Range("G11").Value = 0.5887919840217942
$ws.Range("H11").Value = 0.4844162104220959
$ws.Range("B12").Value = 0.4733464200831813
$ws.Range("C12").Value = 0.6112096612008782
$ws.Range("D12").Value = 0.09205845841638884
$ws.Range("E12").Value = 0.4751725270085909
$ws.Range("F12").Value = 0.4729543549636881
$ws.Range("G12").Value = 0.5416740563608959
$ws.Range("H12").Value = 0.4749351916719513
$ws.Range("B13").Value = 0.5100629818656438
$ws.Range("C13").Value = 0.5945170037428875
$ws.Range("D13").Value = 0.3250309684184554
$ws.Range("E13").Value = 0.4840699376551185
$ws.Range("F13").Value = 0.508822791869401
$ws.Range("G13").Value = 0.5790612702473026
$ws.Range("H13").Value = 0.4838181537396462
